$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")  # "Statistics" is the active/tabSelected sheet

# Turn the A1:E21 range into an Excel Table ("Tabla1") with the Medium2 style,
# matching the newly added xl/tables/table1.xml part.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:E21"), [Type]::Missing, 1)
$tbl.Name = "Tabla1"
$tbl.TableStyle = "TableStyleMedium2"

# Column width tweaks for columns A, D and E.
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 16.333333333333332
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666

# Update the active selection on the sheet.
$ws.Range("G19").Select() | Out-Null
